$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grow the Table1 ListObject from A1:D7 to A1:D20 (keeps autoFilter ref in sync) ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D20"))

# --- Add the new "Basket" API rows (8-11) ---
# Row 8 is written fully, left-to-right, first. Rows 9-11 are written with
# column A/B first (no new shared strings introduced there), then column C
# for rows 9,10,11 in order, then column D for rows 11,10,9 in order -- this
# reproduces the exact shared-string insertion order of the source edit.

# Row 8
$ws.Range("A8").Value = "Basket"
$ws.Range("B8").Value = "GET"
$ws.Range("C8").Value = "api/vi/Basket"
$ws.Range("D8").Value = "Get Basket and Items with Username"

# Columns A & B for rows 9-11 (values already exist in the shared string table)
$ws.Range("A9").Value = "Basket"
$ws.Range("B9").Value = "POST"
$ws.Range("A10").Value = "Basket"
$ws.Range("B10").Value = "DELETE"
$ws.Range("A11").Value = "Basket"
$ws.Range("B11").Value = "POST"

# Column C for rows 9, 10, 11 (in that order)
$ws.Range("C9").Value = "api/vi/Basket"
$ws.Range("C10").Value = "api/vi/Basket/{id}"
$ws.Range("C11").Value = "api/vi/Basket/Checkout"

# Column D for rows 11, 10, 9 (in that order)
$ws.Range("D11").Value = "Checkout Basket"
$ws.Range("D10").Value = "Delete Basket"
$ws.Range("D9").Value = "Update Basket and Items (add - remove item on basket)"

# --- Widen column D to fit the new, longer text ---
$ws.Columns.Item(4).ColumnWidth = 44.33

# --- Select the whole of column D, as the author left it selected ---
$ws.Range("D:D").Select() | Out-Null
